$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 44: F44 becomes a formatted text amount (right under the new F41 total block) ---
$r44 = $ws.Range("F44")
$r44.Value = "4,80,596.00"
$r44.HorizontalAlignment = -4108
$r44.VerticalAlignment = -4160

# --- Row 41: new "Total Amount" row with a SUM formula under the F column ---
$ws.Range("E41").Value = "Total Amount"
$ws.Range("F41").Formula = "=SUM(F5:F40)"

# --- Row 45: F45 becomes a formatted text amount too, keeping its vertical=top and adding center ---
$r45 = $ws.Range("F45")
$r45.Value = "15,00,000.00"
$r45.HorizontalAlignment = -4108
$r45.VerticalAlignment = -4160

# --- Row 46: new "Total Amount" label + grand total text amount ---
$ws.Range("E46").Value = "Total Amount"
$r46 = $ws.Range("F46")
$r46.Value = "19,80,596.00"
$r46.NumberFormat = "#,##0"
$r46.HorizontalAlignment = -4108

# --- Column F widened to fit the new total text ---
$ws.Columns.Item(6).ColumnWidth = 28.25

# --- Sheet view: scroll down to keep the new totals visible, select F47 ---
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F47").Select()
